$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency market data values
$ws.Range("D2").Value = "66.710.44"
$ws.Range("E2").Value = "  +6.29%  "
$ws.Range("D3").Value = "3.589.97"
$ws.Range("E3").Value = "  +3.45%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "415.54"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.36%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "130.25"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.38%  "
$ws.Range("E7").Value = "  +3.51%  "
$ws.Range("D8").Value = "3.581.66"
$ws.Range("E8").Value = "  +3.36%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.999"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.07%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.772"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +6.11%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.174"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +18.46%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000331"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +50.33%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "42.43"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.54%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.88"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.86%  "
$ws.Range("D15").Value = "4.165.53"
$ws.Range("E15").Value = "  +3.56%  "
$ws.Range("E16").Value = "  -0.35%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "20.33"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.04%  "
$ws.Range("D18").Value = "3.630.68"
$ws.Range("E18").Value = "  +5.30%  "
$ws.Range("E19").Value = "  +5.68%  "
$ws.Range("D20").Value = "66.602.08"
$ws.Range("E20").Value = "  +6.13%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.34"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.32%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "448.43"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.35%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "88.90"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.04%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.16"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.50%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "13.07"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.52%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.33"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.49%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.94"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -6.72%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "35.28"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +5.74%  "
$ws.Range("E29").Value = "  +1.27%  "
$ws.Range("E30").Value = "  +3.78%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "12.34"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.89%  "
$ws.Range("E32").Value = "  +4.48%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.38"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.52%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.161"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.30%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "39.94"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.25%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.999"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.06%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "56.67"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.42%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0494"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.52%  "
$ws.Range("D39").Value = "0.0₃0718"
$ws.Range("E39").Value = "  +30.06%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.147"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +9.60%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.997"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.22%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.00"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.66%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "149.30"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.29%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.73"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.27%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.27"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.81%  "
$ws.Range("B46").Value = "TheGraph"
$ws.Range("C46").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.313"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.45%  "
$ws.Range("B47").Value = "NEARProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.32"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.66%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.98"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.55%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.31"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.59%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "15.64"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.59%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "114.32"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.81%  "
